# Auto-generated script applying numeric updates to the Belias_Profits sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3034.6155
$ws.Range("I76").Value = 3032
$ws.Range("K76").Value = 3032
$ws.Range("M76").Value = -2717
$ws.Range("H79").Value = 3034.6155
$ws.Range("I79").Value = 3032
$ws.Range("K79").Value = 3032
$ws.Range("M79").Value = -1940

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 36127.875
$ws.Range("J9").Value = 21799.8
$ws.Range("L9").Value = 21799.8
$ws.Range("N9").Value = -22139.8
$ws.Range("H20").Value = 36127.875
$ws.Range("J20").Value = 21799.8
$ws.Range("L20").Value = 21799.8
$ws.Range("N20").Value = -22339.8
$ws.Range("H63").Value = 3712
$ws.Range("I63").Value = 3156.8572
$ws.Range("J63").Value = 4197.75
$ws.Range("K63").Value = 3156.8572
$ws.Range("L63").Value = 4197.75
$ws.Range("M63").Value = -2470.8572
$ws.Range("N63").Value = -5569.75
$ws.Range("H66").Value = 3712
$ws.Range("I66").Value = 3156.8572
$ws.Range("J66").Value = 4197.75
$ws.Range("K66").Value = 15784.286
$ws.Range("L66").Value = 20988.75
$ws.Range("M66").Value = -12352.286
$ws.Range("N66").Value = -27852.75
$ws.Range("H135").Value = 25572.9
$ws.Range("J135").Value = 25572.9
$ws.Range("L135").Value = 25572.9
$ws.Range("N135").Value = -35712.9

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29169.5
$ws.Range("J55").Value = 29169.5
$ws.Range("L55").Value = 29169.5
$ws.Range("N55").Value = -29715.5
$ws.Range("H70").Value = 93300.164
$ws.Range("J70").Value = 93300.164
$ws.Range("L70").Value = 93300.164
$ws.Range("N70").Value = -93886.164
$ws.Range("H73").Value = 93300.164
$ws.Range("J73").Value = 93300.164
$ws.Range("L73").Value = 93300.164
$ws.Range("N73").Value = -95328.164

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 85073.5
$ws.Range("I2").Value = 5040
$ws.Range("J2").Value = 101080.2
$ws.Range("K2").Value = 30240
$ws.Range("L2").Value = 606481.2
$ws.Range("M2").Value = -30127
$ws.Range("N2").Value = -606707.2
$ws.Range("H5").Value = 955.1429000000001
$ws.Range("I5").Value = 188.90909
$ws.Range("J5").Value = 1798
$ws.Range("K5").Value = 566.72727
$ws.Range("L5").Value = 5394
$ws.Range("M5").Value = -454.72727
$ws.Range("N5").Value = -5618
$ws.Range("H9").Value = 1375.5
$ws.Range("I9").Value = 901
$ws.Range("J9").Value = 1850
$ws.Range("K9").Value = 2703
$ws.Range("L9").Value = 5550
$ws.Range("M9").Value = -2479
$ws.Range("N9").Value = -5998
$ws.Range("H16").Value = 550
$ws.Range("I16").Value = 66.666664
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 199.999992
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = -26.99999199999999
$ws.Range("N16").Value = -6346
$ws.Range("H17").Value = 1579
$ws.Range("I17").Value = 198
$ws.Range("J17").Value = 1751.625
$ws.Range("K17").Value = 594
$ws.Range("L17").Value = 5254.875
$ws.Range("M17").Value = -425
$ws.Range("N17").Value = -5592.875
$ws.Range("H20").Value = 3028.5715
$ws.Range("I20").Value = 800
$ws.Range("J20").Value = 4700
$ws.Range("K20").Value = 2400
$ws.Range("L20").Value = 14100
$ws.Range("M20").Value = -2173
$ws.Range("N20").Value = -14554
$ws.Range("H34").Value = 1760
$ws.Range("I34").Value = 650
$ws.Range("J34").Value = 3980
$ws.Range("K34").Value = 1950
$ws.Range("L34").Value = 11940
$ws.Range("M34").Value = -1866
$ws.Range("N34").Value = -12108
$ws.Range("H39").Value = 1462.1875
$ws.Range("J39").Value = 1764.4445
$ws.Range("L39").Value = 5293.333500000001
$ws.Range("N39").Value = -5881.333500000001
$ws.Range("H55").Value = 682.125
$ws.Range("I55").Value = 364.25
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 1092.75
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = -915.75
$ws.Range("N55").Value = -3354
$ws.Range("H87").Value = 5110
$ws.Range("I87").Value = 4461.6665
$ws.Range("J87").Value = 9000
$ws.Range("K87").Value = 13384.9995
$ws.Range("L87").Value = 27000
$ws.Range("M87").Value = -12136.9995
$ws.Range("N87").Value = -29496
$ws.Range("H90").Value = 5110
$ws.Range("I90").Value = 4461.6665
$ws.Range("J90").Value = 9000
$ws.Range("K90").Value = 40154.9985
$ws.Range("L90").Value = 81000
$ws.Range("M90").Value = -33914.9985
$ws.Range("N90").Value = -93480
$ws.Range("H135").Value = 955.1429000000001
$ws.Range("I135").Value = 188.90909
$ws.Range("J135").Value = 1798
$ws.Range("K135").Value = 1700.18181
$ws.Range("L135").Value = 16182
$ws.Range("M135").Value = 834.81819
$ws.Range("N135").Value = -21252

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 9000
$ws.Range("J20").Value = 9000
$ws.Range("L20").Value = 9000
$ws.Range("N20").Value = -9490
$ws.Range("H24").Value = 9000
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("H57").Value = 19025
$ws.Range("J57").Value = 19025
$ws.Range("L57").Value = 19025
$ws.Range("N57").Value = -20665
$ws.Range("H93").Value = 18066.934
$ws.Range("J93").Value = 18066.934
$ws.Range("L93").Value = 18066.934
$ws.Range("N93").Value = -21810.934
$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524
$ws.Range("M24").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2043.8889
$ws.Range("I16").Value = 2074.5
$ws.Range("J16").Value = 1799
$ws.Range("K16").Value = 2074.5
$ws.Range("L16").Value = 1799
$ws.Range("M16").Value = -1904.5
$ws.Range("N16").Value = -2139
$ws.Range("H22").Value = 774.13336
$ws.Range("I22").Value = 583.3333
$ws.Range("J22").Value = 901.3333
$ws.Range("K22").Value = 583.3333
$ws.Range("L22").Value = 901.3333
$ws.Range("M22").Value = -288.3333
$ws.Range("N22").Value = -1491.3333
$ws.Range("H27").Value = 774.13336
$ws.Range("I27").Value = 583.3333
$ws.Range("J27").Value = 901.3333
$ws.Range("K27").Value = 583.3333
$ws.Range("L27").Value = 901.3333
$ws.Range("M27").Value = -476.3333
$ws.Range("N27").Value = -1115.3333
$ws.Range("H82").Value = 1056.6
$ws.Range("I82").Value = 975.3333
$ws.Range("J82").Value = 1091.4286
$ws.Range("K82").Value = 975.3333
$ws.Range("L82").Value = 1091.4286
$ws.Range("M82").Value = -614.3333
$ws.Range("N82").Value = -1813.4286
$ws.Range("H85").Value = 1056.6
$ws.Range("I85").Value = 975.3333
$ws.Range("J85").Value = 1091.4286
$ws.Range("K85").Value = 975.3333
$ws.Range("L85").Value = 1091.4286
$ws.Range("M85").Value = 272.6667
$ws.Range("N85").Value = -3587.4286
$ws.Range("H88").Value = 29725
$ws.Range("J88").Value = 29725
$ws.Range("L88").Value = 29725
$ws.Range("N88").Value = -30581
$ws.Range("H91").Value = 29725
$ws.Range("J91").Value = 29725
$ws.Range("L91").Value = 29725
$ws.Range("N91").Value = -32689

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 14124.75
$ws.Range("J28").Value = 14124.75
$ws.Range("L28").Value = 14124.75
$ws.Range("N28").Value = -14820.75
$ws.Range("H42").Value = 13000
$ws.Range("J42").Value = 13000
$ws.Range("L42").Value = 13000
$ws.Range("N42").Value = -13756
$ws.Range("H103").Value = 29240.8
$ws.Range("J103").Value = 29240.8
$ws.Range("L103").Value = 29240.8
$ws.Range("N103").Value = -31584.8
